$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 236, shifting the existing rows 236-283
# down to become rows 238-285 (so their data is preserved unchanged).
$ws.Rows("236:237").Insert()

# Fill in the two newly inserted rows with the new weekly price records.

# Row 236: Vega Modelo de Temuco - Zapallo italiano - "Bola 8"
$ws.Range("A236").Value = 10
$ws.Range("B236").Value = "Vega Modelo de Temuco"
$ws.Range("C236").Value = "La Araucanía"
$ws.Range("D236").Value = 44504
$ws.Range("E236").Value = 9
$ws.Range("F236").Value = 100112032
$ws.Range("G236").Value = "Zapallo italiano"
$ws.Range("H236").Value = "Bola 8"
$ws.Range("I236").Value = "Primera"
$ws.Range("J236").Value = 85
$ws.Range("K236").Value = 12000
$ws.Range("L236").Value = 12000
$ws.Range("M236").Value = 12000
$ws.Range("N236").Value = "`$/caja 60 unidades"
$ws.Range("O236").Value = "Región de Arica y Parinacota"
$ws.Range("P236").Value = 200
$ws.Range("Q236").Value = 60
$ws.Range("R236").Value = "Hortaliza"

# Row 237: Vega Modelo de Temuco - Zapallo italiano - "Sin especificar"
$ws.Range("A237").Value = 10
$ws.Range("B237").Value = "Vega Modelo de Temuco"
$ws.Range("C237").Value = "La Araucanía"
$ws.Range("D237").Value = 44504
$ws.Range("E237").Value = 9
$ws.Range("F237").Value = 100112032
$ws.Range("G237").Value = "Zapallo italiano"
$ws.Range("H237").Value = "Sin especificar"
$ws.Range("I237").Value = "Primera"
$ws.Range("J237").Value = 155
$ws.Range("K237").Value = 10000
$ws.Range("L237").Value = 10000
$ws.Range("M237").Value = 10000
$ws.Range("N237").Value = "`$/caja 60 unidades"
$ws.Range("O237").Value = "Región de Arica y Parinacota"
$ws.Range("P237").Value = 167
$ws.Range("Q237").Value = 60
$ws.Range("R237").Value = "Hortaliza"
